# Fix SortExample.xlsx test data:
#  - cells that are blank must not carry a shared-string type (t="s") with
#    no value; they are cleared to true blanks while keeping their fill color.
#  - a few rows had their content in the wrong order (unstable/incorrect sort);
#    those cells are rewritten with the correct value + fill color.
$wb = $excel.ActiveWorkbook

# RGB color (as used by cellXfs fill) for each "style id" used across the sheets,
# converted to the BGR integer Interior.Color expects.
$styleColor = @{
    2 = 9498256   # 90EE90
    3 = 16748574   # 1E90FF
    4 = 13749760   # 00CED1
    5 = 11119017   # A9A9A9
    6 = 8894686   # DEB887
    7 = 6053069   # CD5C5C
    8 = 8034025   # E9967A
    9 = 9639167   # FF1493
}


$ws = $wb.Worksheets.Item("Table")
$ws.Range("B4").Value = ""
$ws.Range("B4").Interior.Color = $styleColor[2]
$ws.Range("C5").Value = ""
$ws.Range("C5").Interior.Color = $styleColor[4]
$ws.Range("F7").Value = ""
$ws.Range("F7").Interior.Color = $styleColor[8]
$ws.Range("A8").Value = ""
$ws.Range("A8").Interior.Color = $styleColor[8]
$ws.Range("H8").Value = ""
$ws.Range("H8").Interior.Color = $styleColor[4]
$ws.Range("H9").Value = ""
$ws.Range("H9").Interior.Color = $styleColor[9]
$ws.Range("C11").Value = ""
$ws.Range("C11").Interior.Color = $styleColor[9]
$ws.Range("G11").Value = ""
$ws.Range("G11").Interior.Color = $styleColor[2]

$ws = $wb.Worksheets.Item("Sort Left to Right")
$ws.Range("E3").Value = ""
$ws.Range("E3").Interior.Color = $styleColor[8]
$ws.Range("A4").Value = ""
$ws.Range("A4").Interior.Color = $styleColor[2]
$ws.Range("B5").Value = ""
$ws.Range("B5").Interior.Color = $styleColor[4]
$ws.Range("H5").Value = ""
$ws.Range("H5").Interior.Color = $styleColor[9]
$ws.Range("A8").Value = "A"
$ws.Range("A8").Interior.Color = $styleColor[5]
$ws.Range("B8").Value = "a"
$ws.Range("B8").Interior.Color = $styleColor[7]
$ws.Range("H8").Value = ""
$ws.Range("H8").Interior.Color = $styleColor[8]
$ws.Range("A9").Value = "A"
$ws.Range("A9").Interior.Color = $styleColor[5]
$ws.Range("B9").Value = "a"
$ws.Range("B9").Interior.Color = $styleColor[7]
$ws.Range("G9").Value = ""
$ws.Range("G9").Interior.Color = $styleColor[2]
$ws.Range("A10").Value = "a"
$ws.Range("A10").Interior.Color = $styleColor[5]
$ws.Range("B10").Value = "A"
$ws.Range("B10").Interior.Color = $styleColor[7]
$ws.Range("D10").Value = ""
$ws.Range("D10").Interior.Color = $styleColor[4]
$ws.Range("F10").Value = ""
$ws.Range("F10").Interior.Color = $styleColor[9]

$ws = $wb.Worksheets.Item("Complex 2")
$ws.Range("B5").Value = ""
$ws.Range("B5").Interior.Color = $styleColor[2]
$ws.Range("F5").Value = ""
$ws.Range("F5").Interior.Color = $styleColor[8]
$ws.Range("C6").Value = ""
$ws.Range("C6").Interior.Color = $styleColor[4]
$ws.Range("F6").Value = "a"
$ws.Range("F6").Interior.Color = $styleColor[6]
$ws.Range("G6").Value = "B"
$ws.Range("G6").Interior.Color = $styleColor[6]
$ws.Range("H6").Value = "A"
$ws.Range("H6").Interior.Color = $styleColor[6]
$ws.Range("F7").Value = "a"
$ws.Range("F7").Interior.Color = $styleColor[7]
$ws.Range("G7").Value = "a"
$ws.Range("G7").Interior.Color = $styleColor[7]
$ws.Range("H7").Value = "A"
$ws.Range("H7").Interior.Color = $styleColor[7]
$ws.Range("A9").Value = ""
$ws.Range("A9").Interior.Color = $styleColor[8]
$ws.Range("H10").Value = ""
$ws.Range("H10").Interior.Color = $styleColor[4]
$ws.Range("G11").Value = ""
$ws.Range("G11").Interior.Color = $styleColor[2]
$ws.Range("C12").Value = ""
$ws.Range("C12").Interior.Color = $styleColor[9]
$ws.Range("H12").Value = ""
$ws.Range("H12").Interior.Color = $styleColor[9]

$ws = $wb.Worksheets.Item("Complex 1")
$ws.Range("B3").Value = ""
$ws.Range("B3").Interior.Color = $styleColor[2]
$ws.Range("H3").Value = ""
$ws.Range("H3").Interior.Color = $styleColor[9]
$ws.Range("C4").Value = ""
$ws.Range("C4").Interior.Color = $styleColor[4]
$ws.Range("H4").Value = ""
$ws.Range("H4").Interior.Color = $styleColor[4]
$ws.Range("F6").Value = ""
$ws.Range("F6").Interior.Color = $styleColor[8]
$ws.Range("A7").Value = ""
$ws.Range("A7").Interior.Color = $styleColor[8]
$ws.Range("C10").Value = ""
$ws.Range("C10").Interior.Color = $styleColor[9]
$ws.Range("G10").Value = ""
$ws.Range("G10").Interior.Color = $styleColor[2]

$ws = $wb.Worksheets.Item("Simple Column")
$ws.Range("A6").Value = ""
$ws.Range("A6").Interior.Color = $styleColor[5]
$ws.Range("A7").Value = ""
$ws.Range("A7").Interior.Color = $styleColor[8]
$ws.Range("D9").Value = ""
$ws.Range("D9").Interior.Color = $styleColor[5]
$ws.Range("D10").Value = ""
$ws.Range("D10").Interior.Color = $styleColor[8]

$ws = $wb.Worksheets.Item("Simple")
$ws.Range("B3").Value = ""
$ws.Range("B3").Interior.Color = $styleColor[2]
$ws.Range("F3").Value = "A"
$ws.Range("F3").Interior.Color = $styleColor[5]
$ws.Range("G3").Value = "A"
$ws.Range("G3").Interior.Color = $styleColor[5]
$ws.Range("H3").Value = "a"
$ws.Range("H3").Interior.Color = $styleColor[5]
$ws.Range("C4").Value = ""
$ws.Range("C4").Interior.Color = $styleColor[4]
$ws.Range("F4").Value = "a"
$ws.Range("F4").Interior.Color = $styleColor[7]
$ws.Range("G4").Value = "a"
$ws.Range("G4").Interior.Color = $styleColor[7]
$ws.Range("H4").Value = "A"
$ws.Range("H4").Interior.Color = $styleColor[7]
$ws.Range("H6").Value = ""
$ws.Range("H6").Interior.Color = $styleColor[4]
$ws.Range("A7").Value = ""
$ws.Range("A7").Interior.Color = $styleColor[8]
$ws.Range("H8").Value = ""
$ws.Range("H8").Interior.Color = $styleColor[9]
$ws.Range("G9").Value = ""
$ws.Range("G9").Interior.Color = $styleColor[2]
$ws.Range("C10").Value = ""
$ws.Range("C10").Interior.Color = $styleColor[9]
$ws.Range("F10").Value = ""
$ws.Range("F10").Interior.Color = $styleColor[8]
